$wb = $excel.ActiveWorkbook
$players = $wb.Worksheets.Item("Players")
$ownerTotals = $wb.Worksheets.Item("OwnerTotals")

# --- Players sheet updates ---
# Row 6
$players.Range("G6").Value = "1:51 - 1st Half"
$players.Range("H6").Value = 2
$players.Range("I6").Value = 4
$players.Range("L6").Value = 1
$players.Range("O6").Value = 1
$players.Range("P6").Value = 11
$players.Range("Q6").Value = 1
$players.Range("R6").Value = 5

# Row 8
$players.Range("G8").Value = "1:51 - 1st Half"

# Row 22
$players.Range("G22").Value = "1:51 - 1st Half"
$players.Range("H22").Value = 15
$players.Range("I22").Value = 10
$players.Range("P22").Value = 15
$players.Range("U22").Value = 4
$players.Range("V22").Value = 4

# Row 26
$players.Range("G26").Value = "1:51 - 1st Half"
$players.Range("H26").Value = 2
$players.Range("J26").Value = 1

# Row 29
$players.Range("G29").Value = "1:51 - 1st Half"
$players.Range("H29").Value = 6
$players.Range("J29").Value = 1
$players.Range("N29").Value = 1
$players.Range("P29").Value = 18
$players.Range("R29").Value = 7

# Row 35
$players.Range("G35").Value = "1:51 - 1st Half"
$players.Range("P35").Value = 11

# Row 39
$players.Range("G39").Value = "1:51 - 1st Half"
$players.Range("H39").Value = 7
$players.Range("I39").Value = 5
$players.Range("P39").Value = 14
$players.Range("U39").Value = 2
$players.Range("V39").Value = 2

# Row 43
$players.Range("G43").Value = "1:51 - 1st Half"
$players.Range("H43").Value = 4
$players.Range("J43").Value = 6

# Row 50
$players.Range("G50").Value = "1:51 - 1st Half"

# Row 52
$players.Range("G52").Value = "1:51 - 1st Half"
$players.Range("O52").Value = 2
$players.Range("P52").Value = 12

# Row 53
$players.Range("G53").Value = "1:51 - 1st Half"
$players.Range("H53").Value = 2
$players.Range("O53").Value = 2
$players.Range("P53").Value = 17
$players.Range("R53").Value = 8

# Row 80
$players.Range("G80").Value = "1:51 - 1st Half"
$players.Range("O80").Value = 3

# Row 84
$players.Range("D84").Value = "London Jemison"
$players.Range("E84").Value = "ALA"
$players.Range("G84").Value = "1:51 - 1st Half"
$players.Range("J84").Value = 2
$players.Range("K84").Value = 0
$players.Range("L84").Value = 0
$players.Range("N84").Value = 0
$players.Range("O84").Value = 1
$players.Range("P84").Value = 8
$players.Range("R84").Value = 3
$players.Range("T84").Value = 1

# Row 85
$players.Range("D85").Value = "Nic Codie"
$players.Range("E85").Value = "TEX"
$players.Range("F85").Value = "TEX@UGA"
$players.Range("G85").Value = "Final"
$players.Range("I85").Value = 2
$players.Range("K85").Value = 1
$players.Range("M85").Value = 2
$players.Range("N85").Value = 1
$players.Range("O85").Value = 0
$players.Range("P85").Value = 14
$players.Range("Q85").Value = 1
$players.Range("R85").Value = 2
$players.Range("T85").Value = 0

# Row 86
$players.Range("D86").Value = "Dylan James"
$players.Range("E86").Value = "UGA"
$players.Range("H86").Value = 4
$players.Range("I86").Value = 0
$players.Range("J86").Value = 3
$players.Range("L86").Value = 1
$players.Range("M86").Value = 1
$players.Range("N86").Value = 0
$players.Range("O86").Value = 4
$players.Range("P86").Value = 11
$players.Range("Q86").Value = 0

# Row 87
$players.Range("D87").Value = "Isaiah Brown"
$players.Range("E87").Value = "FLA"
$players.Range("F87").Value = "FLA@MISS"
$players.Range("I87").Value = 5
$players.Range("K87").Value = 0
$players.Range("L87").Value = 0
$players.Range("M87").Value = 0
$players.Range("N87").Value = 2
$players.Range("O87").Value = 3
$players.Range("P87").Value = 13
$players.Range("Q87").Value = 2
$players.Range("R87").Value = 3
$players.Range("U87").Value = 1
$players.Range("V87").Value = 2

# Row 88
$players.Range("D88").Value = "Jalen Reece"
$players.Range("E88").Value = "LSU"
$players.Range("F88").Value = "ALA@LSU"
$players.Range("G88").Value = "1:51 - 1st Half"
$players.Range("I88").Value = 4
$players.Range("J88").Value = 1
$players.Range("K88").Value = 2
$players.Range("L88").Value = 2
$players.Range("O88").Value = 0
$players.Range("P88").Value = 18
$players.Range("R88").Value = 5
$players.Range("T88").Value = 3
$players.Range("U88").Value = 0
$players.Range("V88").Value = 0

# Row 91
$players.Range("G91").Value = "1:51 - 1st Half"
$players.Range("P91").Value = 7

# Row 102
$players.Range("G102").Value = "1:51 - 1st Half"

# Row 109
$players.Range("G109").Value = "1:51 - 1st Half"

# --- OwnerTotals sheet updates ---
$ownerTotals.Range("B3").Value = 65
$ownerTotals.Range("B4").Value = 53
$ownerTotals.Range("B6").Value = 45
$ownerTotals.Range("B7").Value = 44
$ownerTotals.Range("B8").Value = 25
